$d = $word.ActiveDocument

# Locate the insertion point: right after "El vendedor " and before "también..."
$rng = $d.Content
$found = $rng.Find.Execute("El vendedor ", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
$rng.Collapse(0)
$p1 = $rng.Start

# Move the "_GoBack" bookmark to this (still collapsed) insertion point first.
# Word keeps only a single "_GoBack" bookmark at a time, so adding it here
# automatically removes the old one (which sat at the end of the next, empty
# paragraph) and splits the run at this position.
$d.Bookmarks.Add("_GoBack", $d.Range($p1, $p1))

# Type the new text right at the bookmark; it lands just before the bookmark,
# which then continues to mark the (now shifted) insertion point.
$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertAfter("o Administrador ")

# Finally, split "El vendedor " away from "o Administrador " so they remain two
# distinct runs (matching how Word keeps separately-typed runs separate).
# A short-lived bookmark forces the split; removing it again leaves no trace.
$d.Bookmarks.Add("ZZZ_TempSplit", $d.Range($p1, $p1))
$d.Bookmarks("ZZZ_TempSplit").Delete() | Out-Null
